$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "About" sheet: insert explanatory note about what this variable captures,
#    right under "Key to Table 2 Terms" (row 11), pushing the rest down.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Insert 4 rows at row 12: three for the new note text, one blank spacer
# (matches the old blank row that used to separate this block from the next).
$about.Rows("12:15").Insert()

$about.Range("A12").Value = "This variable captures the capital and O&M costs and energy use of "
$about.Range("A13").Value = "carbon capture and sequestration, for both the power and industry"
$about.Range("A14").Value = "sectors."
$about.Range("A12:A14").Font.Bold = $false

# ---------------------------------------------------------------------------
# 2) Rename the per-ton units to "metric ton CO2e" across the three small
#    CC result sheets (capital cost, O&M cost, energy use), and wrap the
#    now-longer header text in B1.
# ---------------------------------------------------------------------------
$capEqpt = $wb.Worksheets.Item("CC-CCoEtSOToCpY")
$capEqpt.Range("B1").Value = "Capital cost ($/(metric ton CO2e*yr))"
$capEqpt.Range("B1").WrapText = $true
$capEqpt.Rows("1").RowHeight = 30

$omCost = $wb.Worksheets.Item("CC-TOMCpTS")
$omCost.Range("B1").Value = "O&M Cost per Ton ($/metric ton CO2e)"
$omCost.Range("B1").WrapText = $true
$omCost.Rows("1").RowHeight = 30

$energyUse = $wb.Worksheets.Item("CC-EUpTCS")
$energyUse.Range("B1").Value = "Energy Use per Ton Sequestered (BTU/metric ton CO2e)"
$energyUse.Range("B1").WrapText = $true
$energyUse.Rows("1").RowHeight = 30
